$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = 0.261
    $ws.Range("E$row").Value = 0.4370000000000001
    $ws.Range("F$row").Value = -0.0201
    $ws.Range("K$row").Value = 946.7
    $ws.Range("L$row").Value = 0.2825717099961197
    $ws.Range("M$row").Value = 473.1
    $ws.Range("N$row").Value = 0.03813815508387815
    $ws.Range("O$row").Value = 0.4997359247913806
    $ws.Range("P$row").Value = 198
    $ws.Range("Q$row").Value = 0.01596143459439415
    $ws.Range("R$row").Value = 0.2091475652265765
    $ws.Range("S$row").Value = 275.1
    $ws.Range("T$row").Value = 0.5814838300570704
    $ws.Range("U$row").Value = 4388.6
    $ws.Range("V$row").Value = 0.3537795548533241
    $ws.Range("W$row").Value = 0.1323056712413003
    $ws.Range("X$row").Value = 0.0596157356696802
    $ws.Range("Y$row").Value = 0.0726899355716201
    $ws.Range("Z$row").Value = 0.419658290953729
    $ws.Range("AB$row").Value = 0.05301303802385098
    $ws.Range("AC$row").Value = -0.05301303802385098
    $ws.Range("AD$row").Value = 7253.3
    $ws.Range("AF$row").Value = 7253.3
    $ws.Range("AG$row").Value = 2864.7
    $ws.Range("AH$row").Value = 0.3689707094240571
    $ws.Range("AI$row").Value = 0.4758727471936282
    $ws.Range("AJ$row").Value = 0.187608057840415
    $ws.Range("AK$row").Value = 0.2639425070253835
}
